# CAP018_MaintainBooking_TestData.xlsx
# "Change in reporting in Hooks"
#
# A new booking row is inserted right after row 2 (pushing the existing
# rows 3-20 down to 4-21). The new row 3 is a duplicate of the original
# row 2 (SEA/LAX/PRIORITY/2199/10/360/...), and row 2 itself is updated
# in place to reflect a different product/commodity/piece/weight
# combination (GENERAL/NONSCR/13/775) while keeping its agent/shipper/
# consignee/execute/tag values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Duplicate row 2 down into a newly inserted row 3 (shifts old rows
#    3-20 down to 4-21, carrying their formatting/styles with them - this
#    is exactly what Excel's "Copy row > Insert Copied Cells" does).
$ws.Rows(2).Copy()
$ws.Rows(3).Insert()

# 2) Update row 2 in place with the new ProductCode/Commodity/Piece/Weight.
$ws.Range("C2").Value = "GENERAL"
$ws.Range("D2").Value = "NONSCR"
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 775

# 3) Move the active selection to I3 (matches the final saved selection).
$ws.Range("I3").Select()
